# Apply updated cryptocurrency price/volume data to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.340.37"
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").Value = "1.722.57"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'239.62"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "'0.4728"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").Value = "'0.2631"
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("D9").Value = "'0.06232"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").Value = "1.716.55"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").Value = "'0.07071"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "'15.30"
$ws.Range("E12").Value = "  +4.42%  "
$ws.Range("D13").Value = "'0.5932"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "'4.409"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'0.9997"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "26.333.64"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("D19").Value = "'0.000006811"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").Value = "1.936.64"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("D22").Value = "'4.568"
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").Value = "'8.785"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").Value = "'5.332"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'134.97"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").Value = "'15.24"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("D27").Value = "'1.408"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("D28").Value = "'1.764"
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("D29").Value = "'106.96"
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("D30").Value = "'4.021"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").Value = "'3.699"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").Value = "'0.07746"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").Value = "'0.04450"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "'0.9761"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("B37").Value = "Quant"
$ws.Range("C37").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D37").Value = "'115.31"
$ws.Range("E37").Value = "  +17.81%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.9267"
$ws.Range("E38").Value = "  +8.63%  "
$ws.Range("D39").Value = "'2.423"
$ws.Range("E39").Value = "  -6.91%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.909"
$ws.Range("E40").Value = "  +4.94%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'0.01470"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").Value = "'5.345"
$ws.Range("E43").Value = "  +15.21%  "
$ws.Range("D44").Value = "'0.3822"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").Value = "'0.1160"
$ws.Range("E45").Value = "  +5.03%  "
$ws.Range("D46").Value = "'6.260"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "'0.05290"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").Value = "'30.60"
$ws.Range("E48").Value = "  +4.14%  "
$ws.Range("D49").Value = "'7.661"
$ws.Range("E49").Value = "  +4.79%  "
$ws.Range("D50").Value = "'0.3393"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("E51").Value = "  +2.04%  "
